$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E are treated as text so numeric-looking values
# (e.g. "311.00", "0.3640") keep their exact original formatting/precision
# instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '24.391.91'
$ws.Range("E2").Value = '  -1.75%  '
$ws.Range("D3").Value = '1.651.75'
$ws.Range("E3").Value = '  -3.33%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '311.00'
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '0.3640'
$ws.Range("E7").Value = '  -3.00%  '
$ws.Range("D8").Value = '47.01'
$ws.Range("E8").Value = '  -5.41%  '
$ws.Range("D9").Value = '0.3245'
$ws.Range("E9").Value = '  -5.69%  '
$ws.Range("D10").Value = '1.120'
$ws.Range("E10").Value = '  -7.20%  '
$ws.Range("D11").Value = '0.07021'
$ws.Range("E11").Value = '  -6.87%  '
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '5.923'
$ws.Range("E13").Value = '  -5.87%  '
$ws.Range("D14").Value = '19.34'
$ws.Range("E14").Value = '  -8.22%  '
$ws.Range("D15").Value = '6.586'
$ws.Range("E15").Value = '  -6.44%  '
$ws.Range("D16").Value = '1.650.75'
$ws.Range("E16").Value = '  -3.29%  '
$ws.Range("D17").Value = '0.00001041'
$ws.Range("E17").Value = '  -8.33%  '
$ws.Range("D18").Value = '0.06605'
$ws.Range("E18").Value = '  -1.95%  '
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = '77.94'
$ws.Range("E20").Value = '  -7.78%  '
$ws.Range("D21").Value = '5.916'
$ws.Range("E21").Value = '  -7.23%  '
$ws.Range("D22").Value = '15.57'
$ws.Range("E22").Value = '  -10.04%  '
$ws.Range("D23").Value = '12.40'
$ws.Range("E23").Value = '  -6.16%  '
$ws.Range("D24").Value = '24.384.00'
$ws.Range("E24").Value = '  -1.78%  '
$ws.Range("D25").Value = '2.479'
$ws.Range("E25").Value = '  +1.25%  '
$ws.Range("D26").Value = '2.333'
$ws.Range("E26").Value = '  -16.30%  '
$ws.Range("D27").Value = '147.17'
$ws.Range("E27").Value = '  -3.09%  '
$ws.Range("E28").Value = '  -9.04%  '
$ws.Range("D29").Value = '1.835.90'
$ws.Range("E29").Value = '  -3.32%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '123.85'
$ws.Range("E30").Value = '  -6.47%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '1.165'
$ws.Range("E31").Value = '  -5.93%  '
$ws.Range("D32").Value = '4.078'
$ws.Range("E32").Value = '  -4.20%  '
$ws.Range("D33").Value = '5.647'
$ws.Range("E33").Value = '  -18.22%  '
$ws.Range("D34").Value = '0.08460'
$ws.Range("E34").Value = '  -3.64%  '
$ws.Range("D35").Value = '1.659'
$ws.Range("E35").Value = '  -9.40%  '
$ws.Range("D36").Value = '12.26'
$ws.Range("E36").Value = '  -10.71%  '
$ws.Range("D37").Value = '5.168'
$ws.Range("E37").Value = '  -7.88%  '
$ws.Range("D38").Value = '1.242'
$ws.Range("E38").Value = '  -2.84%  '
$ws.Range("D39").Value = '0.06019'
$ws.Range("E39").Value = '  -10.03%  '
$ws.Range("D40").Value = '0.02206'
$ws.Range("E40").Value = '  -8.35%  '
$ws.Range("D41").Value = '0.2062'
$ws.Range("E41").Value = '  -7.82%  '
$ws.Range("D42").Value = '8.124'
$ws.Range("E42").Value = '  -12.75%  '
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("D44").Value = '0.5872'
$ws.Range("E44").Value = '  -8.84%  '
$ws.Range("D45").Value = '3.766'
$ws.Range("E45").Value = '  -1.57%  '
$ws.Range("D46").Value = '12.61'
$ws.Range("E46").Value = '  -9.81%  '
$ws.Range("D47").Value = '0.5601'
$ws.Range("E47").Value = '  -9.11%  '
$ws.Range("D48").Value = '122.17'
$ws.Range("E48").Value = '  -6.16%  '
$ws.Range("D49").Value = '1.940'
$ws.Range("E49").Value = '  -9.26%  '
$ws.Range("D50").Value = '0.06884'
$ws.Range("E50").Value = '  -5.83%  '
$ws.Range("D51").Value = '74.52'
$ws.Range("E51").Value = '  -6.52%  '
